# Fruta / hortaliza, semanal
# Update the weekly data rows (D, L, M, N, O, P, R, S) for the "Higo" (Fig)
# sheet at Mercado Mayorista Lo Valledor de Santiago. The underlying
# logic re-shuffles which market date/quality/volume/price/origin record
# lands on each output row of the daily-to-weekly consolidation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44302; L="Primera"; M=340; N=12000; O=13000; P=12500; R="Provincia de Santiago";  S=1786 },
    @{ Row=3;  D=44349; L="Segunda"; M=70;  N=12000; O=12000; P=12000; R="Región Metropolitana";   S=1714 },
    @{ Row=4;  D=44321; L="Primera"; M=140; N=11000; O=12000; P=11500; R="Región Metropolitana";   S=1643 },
    @{ Row=5;  D=44321; L="Segunda"; M=80;  N=8000;  O=8000;  P=8000;  R="Región Metropolitana";   S=1143 },
    @{ Row=6;  D=44312; L="Primera"; M=50;  N=13000; O=13000; P=13000; R="Región Metropolitana";   S=1857 },
    @{ Row=7;  D=44312; L="Segunda"; M=20;  N=11000; O=11000; P=11000; R="Región Metropolitana";   S=1571 },
    @{ Row=8;  D=44306; L="Primera"; M=50;  N=12000; O=12000; P=12000; R="Región Metropolitana";   S=1714 },
    @{ Row=9;  D=44306; L="Segunda"; M=40;  N=9000;  O=9000;  P=9000;  R="Región Metropolitana";   S=1286 },
    @{ Row=10; D=44315; L="Especial"; M=50; N=14000; O=14000; P=14000; R="Región Metropolitana";   S=2000 },
    @{ Row=11; D=44315; L="Primera"; M=80;  N=12000; O=13000; P=12500; R="Región Metropolitana";   S=1786 },
    @{ Row=12; D=44315; L="Segunda"; M=80;  N=10000; O=11000; P=10500; R="Región Metropolitana";   S=1500 },
    @{ Row=13; D=44322; L="Primera"; M=100; N=11000; O=11000; P=11000; R="Región Metropolitana";   S=1571 },
    @{ Row=14; D=44300; L="Primera"; M=150; N=12000; O=13000; P=12500; R="Provincia de Santiago";  S=1786 },
    @{ Row=15; D=44307; L="Primera"; M=70;  N=14000; O=14000; P=14000; R="Región Metropolitana";   S=2000 },
    @{ Row=16; D=44307; L="Segunda"; M=50;  N=10000; O=10000; P=10000; R="Región Metropolitana";   S=1429 },
    @{ Row=17; D=44316; L="Primera"; M=40;  N=13000; O=13000; P=13000; R="Región Metropolitana";   S=1857 },
    @{ Row=18; D=44316; L="Segunda"; M=50;  N=11000; O=11000; P=11000; R="Región Metropolitana";   S=1571 },
    @{ Row=19; D=44335; L="Primera"; M=80;  N=14000; O=14000; P=14000; R="Región Metropolitana";   S=2000 },
    @{ Row=20; D=44314; L="Primera"; M=20;  N=13000; O=13000; P=13000; R="Región Metropolitana";   S=1857 },
    @{ Row=21; D=44314; L="Segunda"; M=45;  N=11000; O=11000; P=11000; R="Región Metropolitana";   S=1571 },
    @{ Row=22; D=44342; L="Segunda"; M=50;  N=12000; O=12000; P=12000; R="Región Metropolitana";   S=1714 },
    @{ Row=23; D=44344; L="Segunda"; M=50;  N=12000; O=12000; P=12000; R="Región Metropolitana";   S=1714 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 4).Value  = $r.D    # D: Fecha
    $ws.Cells.Item($n, 12).Value = $r.L    # L: Calidad
    $ws.Cells.Item($n, 13).Value = $r.M    # M: Volumen
    $ws.Cells.Item($n, 14).Value = $r.N    # N: Precio minimo
    $ws.Cells.Item($n, 15).Value = $r.O    # O: Precio maximo
    $ws.Cells.Item($n, 16).Value = $r.P    # P: Precio promedio ponderado
    $ws.Cells.Item($n, 18).Value = $r.R    # R: Origen
    $ws.Cells.Item($n, 19).Value = $r.S    # S: Precio $/Kg
}
